# Applies the Sat Oct 19 03:32:45 UTC 2024 cryptos-list refresh:
# updates Price (D) / Volume(1h) (E) figures for the existing rows, and
# corrects rows 40-41 which had swapped Stacks <-> PolygonEcosystemToken.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "`'68.541.12"
$ws.Range("E2").Value = "`'  +0.94%  "

# Row 3
$ws.Range("D3").Value = "`'2.655.89"
$ws.Range("E3").Value = "`'  +1.36%  "

# Row 4
$ws.Range("E4").Value = "`'  +0.04%  "

# Row 5
$ws.Range("D5").Value = "`'600.78"
$ws.Range("E5").Value = "`'  +0.95%  "

# Row 6
$ws.Range("D6").Value = "`'155.33"
$ws.Range("E6").Value = "`'  +1.71%  "

# Row 7
$ws.Range("E7").Value = "`'  +0.00%  "

# Row 8
$ws.Range("E8").Value = "`'  +0.80%  "

# Row 9
$ws.Range("D9").Value = "`'2.655.12"
$ws.Range("E9").Value = "`'  +1.38%  "

# Row 10
$ws.Range("E10").Value = "`'  +9.62%  "

# Row 11
$ws.Range("E11").Value = "`'  -0.21%  "

# Row 12
$ws.Range("D12").Value = "`'5.27"
$ws.Range("E12").Value = "`'  +1.59%  "

# Row 13
$ws.Range("D13").Value = "`'0.357"
$ws.Range("E13").Value = "`'  +2.44%  "

# Row 14
$ws.Range("E14").Value = "`'  +3.08%  "

# Row 15
$ws.Range("E15").Value = "`'  +3.16%  "

# Row 16
$ws.Range("D16").Value = "`'3.130.59"
$ws.Range("E16").Value = "`'  +1.08%  "

# Row 17
$ws.Range("D17").Value = "`'68.422.74"
$ws.Range("E17").Value = "`'  +0.90%  "

# Row 18
$ws.Range("D18").Value = "`'2.663.11"
$ws.Range("E18").Value = "`'  +1.57%  "

# Row 19
$ws.Range("D19").Value = "`'11.54"
$ws.Range("E19").Value = "`'  +2.76%  "

# Row 20
$ws.Range("D20").Value = "`'367.45"
$ws.Range("E20").Value = "`'  -1.31%  "

# Row 21
$ws.Range("D21").Value = "`'7.56"
$ws.Range("E21").Value = "`'  +1.88%  "

# Row 22
$ws.Range("D22").Value = "`'4.45"
$ws.Range("E22").Value = "`'  +5.44%  "

# Row 23
$ws.Range("D23").Value = "`'4.94"
$ws.Range("E23").Value = "`'  +2.50%  "

# Row 24
$ws.Range("E24").Value = "`'  +1.75%  "

# Row 25
$ws.Range("D25").Value = "`'73.87"
$ws.Range("E25").Value = "`'  +1.63%  "

# Row 26
$ws.Range("E26").Value = "`'  +0.16%  "

# Row 27
$ws.Range("E27").Value = "`'  +0.61%  "

# Row 28
$ws.Range("E28").Value = "`'  +5.17%  "

# Row 29
$ws.Range("D29").Value = "`'2.776.09"
$ws.Range("E29").Value = "`'  +0.67%  "

# Row 30
$ws.Range("D30").Value = "`'582.66"
$ws.Range("E30").Value = "`'  -1.22%  "

# Row 31
$ws.Range("D31").Value = "`'1.00"
$ws.Range("E31").Value = "`'  -0.04%  "

# Row 32
$ws.Range("D32").Value = "`'8.22"
$ws.Range("E32").Value = "`'  +5.53%  "

# Row 33
$ws.Range("E33").Value = "`'  +4.88%  "

# Row 34
$ws.Range("E34").Value = "`'  +2.04%  "

# Row 35
$ws.Range("E35").Value = "`'  +5.25%  "

# Row 36
$ws.Range("E36").Value = "`'  +6.04%  "

# Row 37
$ws.Range("E37").Value = "`'  +0.02%  "

# Row 38
$ws.Range("D38").Value = "`'160.73"
$ws.Range("E38").Value = "`'  +1.75%  "

# Row 39
$ws.Range("D39").Value = "`'19.56"
$ws.Range("E39").Value = "`'  +2.32%  "

# Row 40
$ws.Range("B40").Value = "`'PolygonEcosystemToken"
$ws.Range("C40").Value = "`'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "`'0.376"
$ws.Range("E40").Value = "`'  +2.55%  "

# Row 41
$ws.Range("B41").Value = "`'Stacks"
$ws.Range("C41").Value = "`'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "`'1.91"
$ws.Range("E41").Value = "`'  +0.67%  "

# Row 42
$ws.Range("D42").Value = "`'5.45"
$ws.Range("E42").Value = "`'  +3.58%  "

# Row 43
$ws.Range("D43").Value = "`'2.73"
$ws.Range("E43").Value = "`'  +1.59%  "

# Row 44
$ws.Range("E44").Value = "`'  +11.42%  "

# Row 45
$ws.Range("E45").Value = "`'  +3.51%  "

# Row 47
$ws.Range("D47").Value = "`'40.57"
$ws.Range("E47").Value = "`'  +0.42%  "

# Row 48
$ws.Range("D48").Value = "`'158.59"
$ws.Range("E48").Value = "`'  +1.50%  "

# Row 49
$ws.Range("E49").Value = "`'  +3.73%  "

# Row 50
$ws.Range("E50").Value = "`'  +2.72%  "

# Row 51
$ws.Range("E51").Value = "`'  +3.89%  "
